# Updates generated data (attendee counts / event title) for the
# "广州-漫展信息" workbook, mirroring output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 216
$ws.Range("C6").Value = "广州·少女番同人only4.0"
$ws.Range("F7").Value = 537
$ws.Range("F8").Value = 208
$ws.Range("F11").Value = 120
$ws.Range("F12").Value = 574
$ws.Range("F13").Value = 74
$ws.Range("F14").Value = 1727
$ws.Range("F15").Value = 300
$ws.Range("F16").Value = 1786
$ws.Range("F17").Value = 233
$ws.Range("F18").Value = 481
$ws.Range("F19").Value = 31
$ws.Range("F20").Value = 121
$ws.Range("F21").Value = 127

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 35

# ---- Sheet "本地生活" (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5256

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5256
$ws.Range("F6").Value = 216
$ws.Range("C16").Value = "广州·少女番同人only4.0"
$ws.Range("F17").Value = 537
$ws.Range("F18").Value = 208
$ws.Range("F22").Value = 120
$ws.Range("F25").Value = 574
$ws.Range("F26").Value = 74
$ws.Range("F28").Value = 1727
$ws.Range("F29").Value = 300
$ws.Range("F30").Value = 1786
$ws.Range("F31").Value = 35
$ws.Range("F32").Value = 233
$ws.Range("F33").Value = 481
$ws.Range("F34").Value = 31
$ws.Range("F35").Value = 121
$ws.Range("F37").Value = 127

$wb.Save()
